# Phase 3 metadata refresh for the "symptom-impact-cs" CodeSystem workbook:
#   - bump the "Date" property (B8) to the new generation timestamp
#   - populate the previously-blank "Case Sensitive" value (B20) with "true"
#
# Both target cells must remain plain TEXT (shared-string) cells, matching
# the rest of the "Metadata" property sheet, not native booleans/numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Metadata" sheet is the tab-selected/active one

# --- Date (row 8, column B): 2025-10-02T11:12:29+01:00 -> ...T18:31:12+01:00
# This string doesn't look like a pure number/bool, so it is stored as text as-is.
$ws.Cells.Item(8, 2).Value = "2025-10-02T18:31:12+01:00"

# --- Case Sensitive (row 20, column B): blank -> "true"
# A bare Value = "true"/"false" gets auto-coerced to a native boolean by the
# engine's literal-type inference, so force it to text with a leading
# apostrophe (the same "quote prefix" trick Excel itself uses), then copy the
# number format from the neighboring already-text cell (B19, "Copyright"'s
# value) over it so the cell's style matches its plain siblings exactly
# (no left-over quote-prefix flag on the style).
$ws.Cells.Item(20, 2).Value = "'true"
$ws.Cells.Item(19, 2).Copy()
$ws.Cells.Item(20, 2).PasteSpecial(-4122)  # xlPasteFormats
